$d = $word.ActiveDocument
$nbsp = [char]0x00A0

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# --- Release Coordinator team -------------------------------------------------
Replace-Text "Release Coordinator" "Release-coördinator"
Replace-Text "auditor" "controller"
Replace-Text "security" "beveiliging"

# --- HIVE TEAM: DEVELOPMENT ---------------------------------------------------
Replace-Text "$($nbsp)HIVE TEAM: DEVELOPMENT" "HIVE TEAM: ONTWIKKELING"
Replace-Text "Responsible for building SmartCash and supporting applications." "Verantwoordelijk voor het bouwen van SmartCash en bieden van applicatie ondersteuning."
Replace-Text "Creator of the Dash N Drink Soda Machine & SmartCash POS." "Maker van de Dash N Drink Soda Machine & SmartCash POS."
Replace-Text "C++ Software Engineer" "C ++ Software Engineer"
Replace-Text "Developer" "Ontwikkelaar"

# --- HIVE TEAM: OUTREACH 2 -----------------------------------------------------
Replace-Text "$($nbsp)HIVE TEAM: OUTREACH 2" "HIVE TEAM: OUTREACH 2"
Replace-Text "This team focuses on community building, growth, general user acquisition in South America" "Dit team richt zich op gemeenschapsopbouw, groei en het binnen halen van nieuwe gebruikers in Zuid-Amerika"
Replace-Text "Outreach Support" "Outreach ondersteuning"

# --- HIVE TEAM: SUPPORT & WEB --------------------------------------------------
Replace-Text "$($nbsp)HIVE TEAM: SUPPORT$($nbsp)" "HIVE TEAM: ONDERSTEUNING "
Replace-Text "This Hive is responsible for on-boarding & generalized SmartCash support." "Deze Hive is verantwoordelijk voor on-boarding en algemene SmartCash ondersteuning."
Replace-Text "Alex is a jack of all trades who loves Technology, Graphics, Web Design & Infrastructure." "Alex is een manusje-van-alles die houdt van technologie, grafische vormgeving, webdesign en infrastructuur."
Replace-Text "Fiscal Officer" "Fiscaal specialist"
Replace-Text "Support" "Ondersteuning"
Replace-Text "Assistant Coordinator" "Assistent-coördinator"

# --- Closing paragraph ---------------------------------------------------------
Replace-Text "SmartHive will be the lifeblood of the project, which will allow anyone to get involved and submit proposals, helping to generate organic growth at a grassroots level, creating a bottom-up management structure." "SmartHive zal de levensader van het project zijn, waardoor iedereen kan meedoen en voorstellen kan indienen; dit zorgt voor organische groei op het niveau van de gewone gebruiker, waardoor een bottom-up managementstructuur wordt gecreëerd."
